# excluding lactic acid outliers from analysis
# Flip the sign of the GWP_new_frac_reduction (M) and GWP_RIN_frac_reduction (N)
# values for rows 3 (sc1g) and 4 (oc1g), turning negative values into positive ones.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("M3").Value = 0.4474565608894246
$ws.Range("N3").Value = 0.3812195707296989
$ws.Range("M4").Value = 0.05515900976824201
$ws.Range("N4").Value = 0.04665339070962589
